$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" '45.318.30'
Set-TextValue "E2" '  +0.92%  '
Set-TextValue "D3" '2.368.98'
Set-TextValue "E3" '  -0.63%  '
Set-TextValue "E4" '  +0.22%  '
Set-TextValue "D5" '312.84'
Set-TextValue "E5" '  -2.08%  '
Set-TextValue "D6" '107.26'
Set-TextValue "E6" '  -4.08%  '
Set-TextValue "D7" '0.630'
Set-TextValue "E7" '  -0.50%  '
Set-TextValue "E8" '  +0.14%  '
Set-TextValue "D9" '0.610'
Set-TextValue "E9" '  -3.31%  '
Set-TextValue "D10" '40.67'
Set-TextValue "E10" '  -4.64%  '
Set-TextValue "D11" '0.0919'
Set-TextValue "E11" '  -0.96%  '
Set-TextValue "D12" '8.51'
Set-TextValue "E12" '  -2.17%  '
Set-TextValue "E13" '  +0.57%  '
Set-TextValue "D14" '0.980'
Set-TextValue "E14" '  -3.53%  '
Set-TextValue "D15" '2.728.52'
Set-TextValue "E15" '  -0.71%  '
Set-TextValue "D16" '15.35'
Set-TextValue "E16" '  -2.66%  '
Set-TextValue "D17" '2.368.94'
Set-TextValue "E17" '  -1.03%  '
Set-TextValue "D18" '45.362.69'
Set-TextValue "E18" '  +0.95%  '
Set-TextValue "D19" '13.90'
Set-TextValue "E19" '  +6.25%  '
Set-TextValue "D20" '7.27'
Set-TextValue "E20" '  -5.10%  '
Set-TextValue "E21" '  -1.28%  '
Set-TextValue "D22" '73.38'
Set-TextValue "E22" '  -2.35%  '
Set-TextValue "D23" '3.54'
Set-TextValue "E23" '  -2.50%  '
Set-TextValue "D24" '259.74'
Set-TextValue "E24" '  -3.03%  '
Set-TextValue "D25" '2.35'
Set-TextValue "E25" '  +1.13%  '
Set-TextValue "E26" '  +0.21%  '
Set-TextValue "D27" '11.04'
Set-TextValue "E27" '  -1.96%  '
Set-TextValue "D28" '7.28'
Set-TextValue "E28" '  -4.36%  '
Set-TextValue "D29" '2.31'
Set-TextValue "E29" '  -0.76%  '
Set-TextValue "D30" '0.0987'
Set-TextValue "E30" '  +7.70%  '
Set-TextValue "D31" '22.21'
Set-TextValue "E31" '  -2.22%  '
Set-TextValue "D32" '37.11'
Set-TextValue "E32" '  -5.84%  '
Set-TextValue "D33" '167.27'
Set-TextValue "E33" '  -0.80%  '
Set-TextValue "D34" '2.98'
Set-TextValue "E34" '  +1.29%  '
Set-TextValue "D35" '0.131'
Set-TextValue "E35" '  -1.72%  '
Set-TextValue "E36" '  -0.13%  '
Set-TextValue "D37" '4.69'
Set-TextValue "E37" '  -2.22%  '
Set-TextValue "D38" '3.98'
Set-TextValue "E38" '  +1.44%  '
$ws.Range("B39").Value = 'ARBITRUM'
$ws.Range("C39").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue "D39" '1.88'
Set-TextValue "E39" '  +7.42%  '
Set-TextValue "D40" '2.93'
Set-TextValue "E40" '  +0.07%  '
$ws.Range("B41").Value = 'VeChain'
$ws.Range("C41").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue "D41" '0.0355'
Set-TextValue "E41" '  -3.13%  '
Set-TextValue "D42" '97.96'
Set-TextValue "E42" '  -5.51%  '
Set-TextValue "D43" '69.16'
Set-TextValue "E43" '  -3.72%  '
Set-TextValue "D44" '0.228'
Set-TextValue "E44" '  -5.18%  '
$ws.Range("B45").Value = 'FirstDigitalUSD'
$ws.Range("C45").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextValue "D45" '1.00'
Set-TextValue "E45" '  +0.12%  '
$ws.Range("B46").Value = 'Celestia'
$ws.Range("C46").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
Set-TextValue "D46" '12.77'
Set-TextValue "E46" '  -9.86%  '
Set-TextValue "D47" '1.835.23'
Set-TextValue "E47" '  +10.71%  '
Set-TextValue "D48" '84.00'
Set-TextValue "E48" '  +4.69%  '
Set-TextValue "D49" '5.80'
Set-TextValue "E49" '  +4.32%  '
$ws.Range("B50").Value = 'FraxShare'
$ws.Range("C50").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue "D50" '9.32'
Set-TextValue "E50" '  +2.83%  '
$ws.Range("B51").Value = 'Aave'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue "D51" '110.96'
Set-TextValue "E51" '  -7.21%  '
